$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.527.91"
$ws.Range("E2").Value = "  +1.81%  "
$ws.Range("D3").Value = "1.595.64"
$ws.Range("E3").Value = "  +0.78%  "
$ws.Range("E4").Value = "  +0.53%  "
$ws.Range("D5").Value = "'212.05"
$ws.Range("E5").Value = "  +0.03%  "
$ws.Range("D6").Value = "'0.514"
$ws.Range("E6").Value = "  -1.38%  "
$ws.Range("E7").Value = "  +0.55%  "
$ws.Range("D8").Value = "'26.85"
$ws.Range("E8").Value = "  +5.39%  "
$ws.Range("D9").Value = "'43.61"
$ws.Range("E9").Value = "  -4.97%  "
$ws.Range("D11").Value = "'0.0598"
$ws.Range("E11").Value = "  +0.85%  "
$ws.Range("D12").Value = "'0.0909"
$ws.Range("E12").Value = "  +0.98%  "
$ws.Range("D13").Value = "1.823.80"
$ws.Range("E13").Value = "  +0.98%  "
$ws.Range("D14").Value = "1.602.21"
$ws.Range("E14").Value = "  +1.69%  "
$ws.Range("D15").Value = "29.541.40"
$ws.Range("E15").Value = "  +1.91%  "
$ws.Range("E17").Value = "  +1.05%  "
$ws.Range("D18").Value = "'63.86"
$ws.Range("E18").Value = "  +2.61%  "
$ws.Range("D19").Value = "'241.25"
$ws.Range("E19").Value = "  +2.12%  "
$ws.Range("D20").Value = "'7.55"
$ws.Range("E20").Value = "  +1.57%  "
$ws.Range("D21").Value = "0.0₃0692"
$ws.Range("E21").Value = "  -0.10%  "
$ws.Range("E22").Value = "  +0.69%  "
$ws.Range("E23").Value = "  -0.47%  "
$ws.Range("D24").Value = "'9.23"
$ws.Range("E24").Value = "  +0.65%  "
$ws.Range("E25").Value = "  -0.26%  "
$ws.Range("D26").Value = "'154.86"
$ws.Range("E26").Value = "  +1.64%  "
$ws.Range("D27").Value = "'15.36"
$ws.Range("E27").Value = "  +2.03%  "
$ws.Range("E28").Value = "  -0.04%  "
$ws.Range("D29").Value = "'6.39"
$ws.Range("E29").Value = "  +0.97%  "
$ws.Range("E30").Value = "  +0.49%  "
$ws.Range("D31").Value = "'0.0477"
$ws.Range("E31").Value = "  +2.59%  "
$ws.Range("E32").Value = "  +0.03%  "
$ws.Range("D33").Value = "'3.22"
$ws.Range("E33").Value = "  +0.19%  "
$ws.Range("D34").Value = "'3.14"
$ws.Range("E34").Value = "  +3.39%  "
$ws.Range("D35").Value = "1.429.83"
$ws.Range("E35").Value = "  +0.41%  "
$ws.Range("E36").Value = "  +2.13%  "
$ws.Range("D37").Value = "'1.03"
$ws.Range("E37").Value = "  -1.66%  "
$ws.Range("D38").Value = "'2.84"
$ws.Range("E38").Value = "  +2.23%  "
$ws.Range("E39").Value = "  +0.35%  "
$ws.Range("D40").Value = "'0.0166"
$ws.Range("E40").Value = "  +1.54%  "
$ws.Range("D41").Value = "'0.539"
$ws.Range("E41").Value = "  +2.67%  "
$ws.Range("D42").Value = "'1.96"
$ws.Range("E42").Value = "  +0.88%  "
$ws.Range("D43").Value = "'0.0492"
$ws.Range("E43").Value = "  +6.84%  "
$ws.Range("D44").Value = "'54.16"
$ws.Range("E44").Value = "  +6.72%  "
$ws.Range("B45").Value = "ARBITRUM"
$ws.Range("C45").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D45").Value = "'0.799"
$ws.Range("E45").Value = "  +1.63%  "
$ws.Range("B46").Value = "PaxDollar"
$ws.Range("C46").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D46").Value = "'1.00"
$ws.Range("E46").Value = "  +0.48%  "
$ws.Range("D47").Value = "'0.974"
$ws.Range("E47").Value = "  +16.15%  "
$ws.Range("D48").Value = "'65.47"
$ws.Range("E48").Value = "  +1.10%  "
$ws.Range("D49").Value = "'5.33"
$ws.Range("E49").Value = "  +0.33%  "
$ws.Range("D50").Value = "1.735.71"
$ws.Range("E50").Value = "  +0.96%  "
$ws.Range("D51").Value = "'85.90"
$ws.Range("E51").Value = "  +0.40%  "
